$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (A/B widened, C:E set to a uniform width)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19
$ws.Columns.Item(2).ColumnWidth = 44
$ws.Columns.Item(3).ColumnWidth = 11.1667
$ws.Columns.Item(4).ColumnWidth = 11.1667
$ws.Columns.Item(5).ColumnWidth = 11.1667

# ---------------------------------------------------------------------------
# Row 1 - GWPH / GW Pharmaceuticals PLC
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 14
$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = "GWPH "
$ws.Range("B1").Value = "GW Pharmaceuticals PLC"
$ws.Range("C1").Value = 120.2
$ws.Range("D1").Value = 3.7366000000000001
$ws.Range("E1").Value = 2.0739999999999998

# ---------------------------------------------------------------------------
# Row 2 - CRON / Cronos Group Inc
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14
$ws.Range("A2").Value = "CRON"
$ws.Range("B2").Value = "Cronos Group Inc"
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 2.0935999999999999
$ws.Range("E2").Value = 1.8069999999999999

# ---------------------------------------------------------------------------
# Row 3 - TLRY / Tilray Inc
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 14
$ws.Range("A3").Value = "TLRY"
$ws.Range("B3").Value = "Tilray Inc"
$ws.Range("C3").Value = 7.95
$ws.Range("D3").Value = "992.65M"
$ws.Range("E3").Value = 1.262

# ---------------------------------------------------------------------------
# Trailing blank rows (row 6 intentionally left untouched/absent)
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Style = "Normal"
$ws.Rows.Item(4).RowHeight = 14
$ws.Rows.Item(5).Style = "Normal"
$ws.Rows.Item(5).RowHeight = 14
$ws.Rows.Item(7).Style = "Normal"
$ws.Rows.Item(7).RowHeight = 14
$ws.Rows.Item(8).Style = "Normal"
$ws.Rows.Item(8).RowHeight = 14
$ws.Rows.Item(9).Style = "Normal"
$ws.Rows.Item(9).RowHeight = 14

# ---------------------------------------------------------------------------
# Drop the now-unused named cell styles (Hyperlink / Followed Hyperlink) that
# were only referenced by the old A1 formatting.
# ---------------------------------------------------------------------------
$wb.Styles.Item("Followed Hyperlink").Delete()
$wb.Styles.Item("Hyperlink").Delete()

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$ws.Range("D7").Select()
